$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 8 gseed/sseed values (2100 -> 2101)
$ws.Range("C8").Value = 2101
$ws.Range("D8").Value = 2101

# Add new "type" column (J) and "kernel" column (K) headers
$ws.Range("J6").Value = "type"
$ws.Range("J7").Value = "AS"
$ws.Range("J8").Value = "AS"
$ws.Range("J9").Value = "AS"
$ws.Range("J10").Value = "AS"

$ws.Range("K6").Value = "kernel"
$ws.Range("K7").Value = "Long"
$ws.Range("K8").Value = "Long"
$ws.Range("K9").Value = "Hall"
$ws.Range("K10").Value = "Hall"

# Add new rows 9 and 10 with data
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 3001
$ws.Range("D9").Value = 3001
$ws.Range("E9").Value = 75
$ws.Range("F9").Value = 16
$ws.Range("G9").Value = 24
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = "n"

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 3101
$ws.Range("D10").Value = 3101
$ws.Range("E10").Value = 75
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = 24
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = "y"

# Update selection to match final state
$ws.Range("I12").Select()
